$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the date serial values in column F (rows 2-7) forward by 19 days.
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $cell.Value2 + 19
}
